$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 : 2022-05-05, 15:20 -> 16:30 (Documentation de projet / Ajout des maquettes...) ---
$ws.Range("A19").Value = 44686
$ws.Range("B19").Value = 0.63888888888888895
$ws.Range("C19").Value = 0.6875

# --- Row 20 : 2022-05-05, 16:30 -> 16:55 (Discussion avec chef de projet / Le MCD n'étais...) ---
$ws.Range("A20").Value = 44686
$ws.Range("B20").Value = 0.6875
$ws.Range("C20").Value = 0.70486111111111116

# --- Row 21 : 2022-05-06, 08:00 -> 09:00 (Correction du MCD) ---
$ws.Range("A21").Value = 44687
$ws.Range("B21").Value = 0.33333333333333331
$ws.Range("C21").Value = 0.375

# --- Row 22 : 2022-05-06, 09:00 -> 09:30 ---
$ws.Range("A22").Value = 44687
$ws.Range("B22").Value = 0.375
$ws.Range("C22").Value = 0.39583333333333331

# --- Row 23 / 24 : date only ---
$ws.Range("A23").Value = 44687
$ws.Range("A24").Value = 44687

# Text entries - order matters to reproduce the exact shared-string table
# (E19 reuses an existing string, the rest create new ones in this sequence).
$ws.Range("E19").Value = "Documentation de projet"
$ws.Range("E20").Value = "Discussion avec chef de projet"
$ws.Range("F20").Value = "Le MCD n'étais pas encore au points do"
$ws.Range("F19").Value = "Ajout des maquettes et modification de la stratégie de tests"
$ws.Range("E21").Value = "Correction du MCD"

$excel.Calculate()

$ws.Range("F25").Select()
